$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph ("Students are to render ...") was split into three
#    runs ("...Students ", "are", " then to email ... the following:").
#    Re-merge them into a single run with the identical, already-correct
#    combined text (a same-text Find/Replace across the run boundaries
#    collapses the runs into one).
# ---------------------------------------------------------------------------
$combined = "Students are to render the 3 images defined in the design brief. " +
    "Students are then to email their lecturer the rendered images, which the " +
    "lecturer will then provide feedback on the images. Students are then to " +
    "render the 3 images again, implementing the feedback received. Students " +
    "are then to submit the following:"

$d.Content.Find.Execute($combined, $true, $false, $false, $false, $false,
    $true, 1, $false, $combined, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Two list-item paragraphs ("A document specifying" and "what feedback
#    was received") each have their text split across two runs with a
#    <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/> pair
#    wrapped around the second run (leftover Word grammar-check markup).
#    Re-serializing the paragraph's own Range.WordOpenXML and feeding it
#    back in via InsertXML merges the runs into one and drops the
#    now-orphaned proofErr markers, while fully preserving the paragraph's
#    own formatting (pPr/numbering/paraId/etc. come along unchanged).
# ---------------------------------------------------------------------------
function Clean-Paragraph($para) {
    $r = $para.Range
    $openXml = $r.WordOpenXML

    $startTag = $openXml.IndexOf("<w:p ")
    if ($startTag -lt 0) { $startTag = $openXml.IndexOf("<w:p>") }
    $endTag = $openXml.IndexOf("</w:p>") + 6
    $paraXml = $openXml.Substring($startTag, $endTag - $startTag)

    $wrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
        $paraXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($wrapped)
}

$targets = @("A document specifying", "what feedback was received")

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    foreach ($t in $targets) {
        if ($text -eq $t) {
            Clean-Paragraph $para
            break
        }
    }
}
